# The commit removes one of the scraped "problem" URLs (the
# findonlinejobsopen.space test row) from the scanner test fixture and
# updates the sheet's font from Arial to Calibri. Deleting the whole row
# shifts every row below it up by one and the engine automatically
# reindexes/prunes the now-unused shared string, which matches the
# uniqueCount 86 -> 85 change seen in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row that holds the "findonlinejobsopen.space" URL and remove it
# entirely (rather than hard-coding row 6) so the script is resilient to
# any prior shifting.
$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count
$targetRow = 0
for ($r = 1; $r -le $lastRow; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($val -like "https://findonlinejobsopen.space*") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows.Item($targetRow).Delete() | Out-Null
}

# Switch the sheet's font from Arial to Calibri (the "Normal" style drives
# every cell here since none carry a per-cell override).
$wb.Styles.Item("Normal").Font.Name = "Calibri"

# Move the active selection to E10, as recorded in the saved view state.
$ws.Range("E10").Select() | Out-Null
